$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images_eeg/Sphere_CW-3.75_BG-grey_stim-white.png",
    "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images_eeg/Sphere_Ref_BG-grey_stim-white.png",
    "./images_eeg/Sphere_CW-3.75_BG-grey_stim-white.png",
    "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images_eeg/Sphere_Ref_BG-grey_stim-white.png",
    "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images_eeg/Sphere_CW-3.75_BG-grey_stim-white.png",
    "./images_eeg/Sphere_Ref_BG-grey_stim-white.png",
    "./images_eeg/Sphere_CW-3.75_BG-grey_stim-white.png",
    "./images_eeg/Sphere_CW-3.75_BG-grey_stim-white.png",
    "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images_eeg/Sphere_Ref_BG-grey_stim-white.png",
    "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-white.png",
    "./images_eeg/Sphere_CW-3.75_BG-grey_stim-white.png",
    "./images_eeg/Sphere_CW-3.75_BG-grey_stim-white.png"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
